$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the runs/balls/fours stats between row 2 and row 3,
# keeping values stored as text (matching original t="str" cells).
$textRange = $ws.Range("C2:E3")
$textRange.NumberFormat = "@"

$ws.Range("C2").Value = "4"
$ws.Range("D2").Value = "2"
$ws.Range("E2").Value = "1"

$ws.Range("C3").Value = "0"
$ws.Range("D3").Value = "1"
$ws.Range("E3").Value = "0"
